$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers in H1:J1 from ft_* to ts_*
$ws.Range("H1").Value = "ts_hombres"
$ws.Range("I1").Value = "ts_mujeres"
$ws.Range("J1").Value = "ts_total"

# Append new row 19 with 2019 data (anuario 2019 para OOSS y huelgas)
$ws.Range("A19").Value = 2019
$ws.Range("B19").Value = 689228
$ws.Range("C19").Value = 503876
$ws.Range("D19").Value = 1193104
$ws.Range("E19").Value = 4446632.258064516
$ws.Range("F19").Value = 3072414.634146342
$ws.Range("G19").Value = 7503798.742138364
$ws.Range("H19").Value = 15.5
$ws.Range("I19").Value = 16.4
$ws.Range("J19").Value = 15.9
